# Insert a new data row at row 245 (pushes existing rows 245-257 down to
# 246-258), then populate the newly inserted row with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(245).Insert()

$ws.Range("A245").Value = 4
$ws.Range("B245").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C245").Value = "Los Lagos"
$ws.Range("D245").Value = 44753
$ws.Range("E245").Value = 10
$ws.Range("F245").Value = 100112032
$ws.Range("G245").Value = "Zapallo italiano"
$ws.Range("H245").Value = "Sin especificar"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 70
$ws.Range("K245").Value = 15000
$ws.Range("L245").Value = 15000
$ws.Range("M245").Value = 15000
$ws.Range("N245").Value = "$/caja 50 unidades"
$ws.Range("O245").Value = "Región de Arica y Parinacota"
$ws.Range("P245").Value = 300
$ws.Range("Q245").Value = 50
$ws.Range("R245").Value = "Hortaliza"
